# Sheet-exporting style edit:
#   - the original worksheet is renamed "TRY" (kept as-is, same sheetId/rId)
#   - four new worksheets are appended, each populated with a view of the
#     same small "people" dataset (full table, a row slice, a filtered
#     view, and a column/row "loc" slice)

$wb = $excel.ActiveWorkbook

function Set-StandardMargins($ws) {
    # match the 0.75in/1in/0.5in margins used throughout the workbook
    $ws.PageSetup.LeftMargin = 54
    $ws.PageSetup.RightMargin = 54
    $ws.PageSetup.TopMargin = 72
    $ws.PageSetup.BottomMargin = 72
    $ws.PageSetup.HeaderMargin = 36
    $ws.PageSetup.FooterMargin = 36
}

function Write-Table($ws, $headers, $data) {
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
    }
    for ($r = 0; $r -lt $data.Length; $r++) {
        $row = $data[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }
}

# --- sheet 1 (existing): just rename to "TRY" ------------------------------
$wsTRY = $wb.Worksheets.Item(1)
$wsTRY.Name = "TRY"
$headerStyleSource = $wsTRY.Range("A1")

# --- full dataset ------------------------------------------------------------
$headers = @("Name", "Age", "Occupation", "Score")
$people = @(
    @("Alice",   24, "Data Analyst", 88),
    @("Bob",     30, "Engineer",     92),
    @("Charlie", 22, "Teacher",      79),
    @("Diana",   28, "Designer",     85),
    @("Evan",    35, "Manager",      90)
)

$wsData = $wb.Worksheets.Add($null, $wsTRY)
$wsData.Name = "Sheet1"
Write-Table $wsData $headers $people
$headerStyleSource.Copy()
$wsData.Range("A1:D1").PasteSpecial(-4122)
Set-StandardMargins $wsData

# --- first 3 rows slice -------------------------------------------------------
$wsSlice = $wb.Worksheets.Add($null, $wsData)
$wsSlice.Name = "Slice_Rows_First3"
Write-Table $wsSlice $headers $people[0..2]
$headerStyleSource.Copy()
$wsSlice.Range("A1:D1").PasteSpecial(-4122)
Set-StandardMargins $wsSlice

# --- rows filtered on Age > 25 ------------------------------------------------
$wsFilt = $wb.Worksheets.Add($null, $wsSlice)
$wsFilt.Name = "Filtered_Age>25"
$filtered = @($people[1], $people[3], $people[4])   # Bob, Diana, Evan
Write-Table $wsFilt $headers $filtered
$headerStyleSource.Copy()
$wsFilt.Range("A1:D1").PasteSpecial(-4122)
Set-StandardMargins $wsFilt

# --- .loc style slice: Name/Score columns, rows 2-4 ---------------------------
$wsLoc = $wb.Worksheets.Add($null, $wsFilt)
$wsLoc.Name = "Loc_Slice"
$locHeaders = @("Name", "Score")
$locRows = @(
    @("Bob",     92),
    @("Charlie", 79),
    @("Diana",   85)
)
Write-Table $wsLoc $locHeaders $locRows
$headerStyleSource.Copy()
$wsLoc.Range("A1:B1").PasteSpecial(-4122)
Set-StandardMargins $wsLoc

# keep the originally active sheet/tab selected
$wsTRY.Activate()
